$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2273.75
$ws.Range("J17").Value = 2273.75
$ws.Range("L17").Value = 6821.25
$ws.Range("N17").Value = -7157.25
# Row 18
$ws.Range("H18").Value = 8061.375
$ws.Range("I18").Value = 9184.429
$ws.Range("K18").Value = 9184.429
$ws.Range("M18").Value = -8900.429
# Row 40
$ws.Range("H40").Value = 3573.0
$ws.Range("J40").Value = 4515.1665
$ws.Range("L40").Value = 4515.1665
$ws.Range("N40").Value = -4865.1665
# Row 43
$ws.Range("H43").Value = 12647.125
$ws.Range("J43").Value = 12647.125
$ws.Range("L43").Value = 12647.125
$ws.Range("N43").Value = -12785.125
# Row 53
$ws.Range("H53").Value = 5536.625
$ws.Range("I53").Value = 7115.5
$ws.Range("K53").Value = 7115.5
$ws.Range("M53").Value = -6478.5
# Row 55
$ws.Range("H55").Value = 233.66667
$ws.Range("I55").Value = 100.25
$ws.Range("K55").Value = 100.25
$ws.Range("M55").Value = 113.75
# Row 76
$ws.Range("H76").Value = 5964.923
$ws.Range("I76").Value = 0.0
$ws.Range("K76").Value = 0.0
$ws.Range("M76").ClearContents()
# Row 79
$ws.Range("H79").Value = 5964.923
$ws.Range("I79").Value = 0.0
$ws.Range("K79").Value = 0.0
$ws.Range("M79").ClearContents()
# Row 88
$ws.Range("H88").Value = 874.0
$ws.Range("J88").Value = 867.5
$ws.Range("L88").Value = 867.5
$ws.Range("N88").Value = -1679.5
# Row 91
$ws.Range("H91").Value = 874.0
$ws.Range("J91").Value = 867.5
$ws.Range("L91").Value = 867.5
$ws.Range("N91").Value = -3675.5
# Row 93
$ws.Range("H93").Value = 601000000.0
$ws.Range("J93").Value = 601000000.0
$ws.Range("L93").Value = 601000000.0
$ws.Range("N93").Value = -601004992.0
# Row 107
$ws.Range("H107").Value = 10617.692
$ws.Range("I107").Value = 15518.375
$ws.Range("K107").Value = 15518.375
$ws.Range("M107").Value = -13598.375
# Row 112
$ws.Range("H112").Value = 2595.3333
$ws.Range("I112").Value = 1346.5
$ws.Range("J112").Value = 2952.1428
$ws.Range("K112").Value = 4039.5
$ws.Range("L112").Value = 8856.4284
$ws.Range("M112").Value = -2931.5
$ws.Range("N112").Value = -11072.4284
# Row 125
$ws.Range("H125").Value = 5899.2144
$ws.Range("I125").Value = 10599.0
$ws.Range("J125").Value = 4019.3
$ws.Range("K125").Value = 95391.0
$ws.Range("L125").Value = 36173.7
$ws.Range("M125").Value = -92931.0
$ws.Range("N125").Value = -41093.7

$ws = $wb.Worksheets.Item("ARM")
# Row 88
$ws.Range("H88").Value = 500002500.0
$ws.Range("I88").Value = 0.0
$ws.Range("K88").Value = 0.0
$ws.Range("M88").ClearContents()
# Row 91
$ws.Range("H91").Value = 500002500.0
$ws.Range("I91").Value = 0.0
$ws.Range("K91").Value = 0.0
$ws.Range("M91").ClearContents()
# Row 97
$ws.Range("H97").Value = 7696722.5
$ws.Range("I97").Value = 6717.3125
$ws.Range("K97").Value = 6717.3125
$ws.Range("M97").Value = -6221.3125
# Row 122
$ws.Range("H122").Value = 1261619.8
$ws.Range("I122").Value = 4175.84
$ws.Range("J122").Value = 4405229.5
$ws.Range("K122").Value = 12527.52
$ws.Range("L122").Value = 13215688.5
$ws.Range("M122").Value = -10077.52
$ws.Range("N122").Value = -13220588.5

$ws = $wb.Worksheets.Item("BSM")
# Row 75
$ws.Range("H75").Value = 70034.8
$ws.Range("I75").Value = 52391.332
$ws.Range("J75").Value = 96500.0
$ws.Range("K75").Value = 52391.332
$ws.Range("L75").Value = 96500.0
$ws.Range("M75").Value = -51455.332
$ws.Range("N75").Value = -98372.0
# Row 78
$ws.Range("H78").Value = 70034.8
$ws.Range("I78").Value = 52391.332
$ws.Range("J78").Value = 96500.0
$ws.Range("K78").Value = 157173.996
$ws.Range("L78").Value = 289500.0
$ws.Range("M78").Value = -152493.996
$ws.Range("N78").Value = -298860.0
# Row 86
$ws.Range("H86").Value = 7958.2144
$ws.Range("J86").Value = 2403.5715
$ws.Range("L86").Value = 2403.5715
$ws.Range("N86").Value = -4649.5715
# Row 89
$ws.Range("H89").Value = 7958.2144
$ws.Range("J89").Value = 2403.5715
$ws.Range("L89").Value = 12017.8575
$ws.Range("N89").Value = -23249.8575
# Row 99
$ws.Range("H99").Value = 16105.037
$ws.Range("I99").Value = 19151.55
$ws.Range("K99").Value = 19151.55
$ws.Range("M99").Value = -17653.55
# Row 102
$ws.Range("H102").Value = 11666.0
$ws.Range("I102").Value = 11666.0
$ws.Range("K102").Value = 11666.0
$ws.Range("M102").Value = -8421.0

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1505.8334
$ws.Range("I16").Value = 1262.75
$ws.Range("J16").Value = 1992.0
$ws.Range("K16").Value = 1262.75
$ws.Range("L16").Value = 1992.0
$ws.Range("M16").Value = -975.75
$ws.Range("N16").Value = -2566.0
# Row 113
$ws.Range("H113").Value = 1505.8334
$ws.Range("I113").Value = 1262.75
$ws.Range("J113").Value = 1992.0
$ws.Range("K113").Value = 1262.75
$ws.Range("L113").Value = 1992.0
$ws.Range("M113").Value = 907.25
$ws.Range("N113").Value = -6332.0
# Row 134
$ws.Range("H134").Value = 2851.5454
$ws.Range("I134").Value = 2968.6667
$ws.Range("J134").Value = 2324.5
$ws.Range("K134").Value = 8906.000100000001
$ws.Range("L134").Value = 6973.5
$ws.Range("M134").Value = -6371.000100000001
$ws.Range("N134").Value = -12043.5

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 77.375
$ws.Range("I2").Value = 85.28571
$ws.Range("J2").Value = 22.0
$ws.Range("K2").Value = 511.71426
$ws.Range("L2").Value = 132.0
$ws.Range("M2").Value = -398.71426
$ws.Range("N2").Value = -358.0
# Row 5
$ws.Range("H5").Value = 435774.56
$ws.Range("I5").Value = 422.4
$ws.Range("K5").Value = 1267.2
$ws.Range("M5").Value = -1155.2
# Row 103
$ws.Range("H103").Value = 6574.0586
$ws.Range("J103").Value = 5264.5454
$ws.Range("L103").Value = 15793.6362
$ws.Range("N103").Value = -17551.6362
# Row 107
$ws.Range("H107").Value = 1069.76
$ws.Range("J107").Value = 1335.7778
$ws.Range("L107").Value = 4007.3334
$ws.Range("N107").Value = -7847.3334
# Row 135
$ws.Range("H135").Value = 435774.56
$ws.Range("I135").Value = 422.4
$ws.Range("K135").Value = 3801.6
$ws.Range("M135").Value = -1266.6
# Row 139
$ws.Range("H139").Value = 2729399.2
$ws.Range("I139").Value = 4286770.5
$ws.Range("K139").Value = 12860311.5
$ws.Range("M139").Value = -12855171.5

$ws = $wb.Worksheets.Item("GSM")
# Row 10
$ws.Range("H10").Value = 7999.0
$ws.Range("I10").Value = 7999.0
$ws.Range("J10").Value = 0.0
$ws.Range("K10").Value = 7999.0
$ws.Range("L10").Value = 0.0
$ws.Range("M10").Value = -7830.0
$ws.Range("N10").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 27502.545
$ws.Range("I40").Value = 30639.2
$ws.Range("K40").Value = 30639.2
$ws.Range("M40").Value = -30503.2
# Row 55
$ws.Range("H55").Value = 997.25
$ws.Range("J55").Value = 1619.6364
$ws.Range("L55").Value = 1619.6364
$ws.Range("N55").Value = -1965.6364
# Row 100
$ws.Range("H100").Value = 5713.857
$ws.Range("I100").Value = 2499.75
$ws.Range("J100").Value = 9999.333
$ws.Range("K100").Value = 2499.75
$ws.Range("L100").Value = 9999.333
$ws.Range("M100").Value = -1958.75
$ws.Range("N100").Value = -11081.333
# Row 132
$ws.Range("H132").Value = 483970.1
$ws.Range("I132").Value = 598983.0
$ws.Range("K132").Value = 1796949.0
$ws.Range("M132").Value = -1794419.0
# Row 136
$ws.Range("H136").Value = 3738.8235
$ws.Range("I136").Value = 2464.5454
$ws.Range("J136").Value = 6075.0
$ws.Range("K136").Value = 7393.6362
$ws.Range("L136").Value = 18225.0
$ws.Range("M136").Value = -4843.6362
$ws.Range("N136").Value = -23325.0

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 10401.083
$ws.Range("J81").Value = 3310.6667
$ws.Range("L81").Value = 6621.3334
$ws.Range("N81").Value = -8743.3334
# Row 84
$ws.Range("H84").Value = 10401.083
$ws.Range("J84").Value = 3310.6667
$ws.Range("L84").Value = 33106.667
$ws.Range("N84").Value = -43714.667
# Row 122
$ws.Range("H122").Value = 4940.6895
$ws.Range("I122").Value = 1775.2858
$ws.Range("K122").Value = 5325.857400000001
$ws.Range("M122").Value = -2875.857400000001
# Row 132
$ws.Range("H132").Value = 11746.975
$ws.Range("I132").Value = 12552.277
$ws.Range("J132").Value = 4499.25
$ws.Range("K132").Value = 37656.831
$ws.Range("L132").Value = 13497.75
$ws.Range("M132").Value = -35126.831
$ws.Range("N132").Value = -18557.75
# Row 136
$ws.Range("H136").Value = 596554.4
$ws.Range("I136").Value = 859134.4
$ws.Range("J136").Value = 5749.25
$ws.Range("K136").Value = 2577403.2
$ws.Range("L136").Value = 17247.75
$ws.Range("M136").Value = -2574853.2
$ws.Range("N136").Value = -22347.75
